# Image folder cleanup: several stimulus image files had their text
# labels removed, which changed which generated filename corresponds to
# which set of (conceptual/perceptual/typicality/n/p_*/r_*) statistics.
# Re-point each affected row's stimulus column (L) at its new filename
# and bring along the matching stats columns (M:V).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = "stimuli/img_pey7u.png"
$ws.Range("M3").Value = 30.34883720930232
$ws.Range("N3").Value = 20.34883720930232
$ws.Range("O3").Value = 25.34883720930232
$ws.Range("P3").Value = 43
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 2

$ws.Range("L5").Value = "stimuli/img_5tr4v.png"
$ws.Range("M5").Value = 56.86046511627907
$ws.Range("N5").Value = 39.3953488372093
$ws.Range("O5").Value = 48.12790697674419
$ws.Range("P5").Value = 43
$ws.Range("Q5").Value = 4
$ws.Range("R5").Value = 4
$ws.Range("S5").Value = 4
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = 4
$ws.Range("V5").Value = 4

$ws.Range("L6").Value = "stimuli/img_5nlnv.png"
$ws.Range("M6").Value = 86.1219512195122
$ws.Range("N6").Value = 69.1951219512195
$ws.Range("O6").Value = 77.65853658536585
$ws.Range("P6").Value = 41
$ws.Range("Q6").Value = 9
$ws.Range("R6").Value = 9
$ws.Range("S6").Value = 9
$ws.Range("T6").Value = 9
$ws.Range("U6").Value = 9
$ws.Range("V6").Value = 9

$ws.Range("L8").Value = "stimuli/img_9oofc.png"
$ws.Range("M8").Value = 82.47619047619048
$ws.Range("N8").Value = 65.5
$ws.Range("O8").Value = 73.98809523809524
$ws.Range("P8").Value = 42
$ws.Range("Q8").Value = 8
$ws.Range("R8").Value = 8
$ws.Range("S8").Value = 8
$ws.Range("T8").Value = 8
$ws.Range("U8").Value = 8
$ws.Range("V8").Value = 8

$ws.Range("L9").Value = "stimuli/img_il020.png"
$ws.Range("M9").Value = 18.85416666666667
$ws.Range("N9").Value = 16.16666666666667
$ws.Range("O9").Value = 17.51041666666667
$ws.Range("P9").Value = 48
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 1
$ws.Range("S9").Value = 1
$ws.Range("T9").Value = 1
$ws.Range("U9").Value = 1
$ws.Range("V9").Value = 1

$ws.Range("L15").Value = "stimuli/img_x9w7o.png"
$ws.Range("M15").Value = 92.38888888888889
$ws.Range("N15").Value = 72.94444444444444
$ws.Range("O15").Value = 82.66666666666666
$ws.Range("P15").Value = 36
$ws.Range("Q15").Value = 10
$ws.Range("R15").Value = 10
$ws.Range("S15").Value = 10
$ws.Range("T15").Value = 10
$ws.Range("U15").Value = 10
$ws.Range("V15").Value = 10

$ws.Range("L16").Value = "stimuli/img_qz292.png"
$ws.Range("M16").Value = 78.26666666666667
$ws.Range("N16").Value = 59.13333333333333
$ws.Range("O16").Value = 68.7
$ws.Range("P16").Value = 45
$ws.Range("Q16").Value = 7
$ws.Range("R16").Value = 7
$ws.Range("S16").Value = 7
$ws.Range("T16").Value = 7
$ws.Range("U16").Value = 7
$ws.Range("V16").Value = 7

$ws.Range("L21").Value = "stimuli/img_s2zoe.png"
$ws.Range("M21").Value = 64.71428571428571
$ws.Range("N21").Value = 44.90476190476191
$ws.Range("O21").Value = 54.80952380952381
$ws.Range("P21").Value = 42
$ws.Range("Q21").Value = 5
$ws.Range("R21").Value = 5
$ws.Range("S21").Value = 5
$ws.Range("T21").Value = 5
$ws.Range("U21").Value = 5
$ws.Range("V21").Value = 5

$ws.Range("L23").Value = "stimuli/img_jpjeg.png"
$ws.Range("M23").Value = 90.90697674418605
$ws.Range("N23").Value = 74.3953488372093
$ws.Range("O23").Value = 82.65116279069767
$ws.Range("P23").Value = 43
$ws.Range("Q23").Value = 10
$ws.Range("R23").Value = 10
$ws.Range("S23").Value = 10
$ws.Range("T23").Value = 10
$ws.Range("U23").Value = 10
$ws.Range("V23").Value = 10

$ws.Range("L25").Value = "stimuli/img_rru0v.png"
$ws.Range("M25").Value = 56.45238095238095
$ws.Range("N25").Value = 39.42857142857143
$ws.Range("O25").Value = 47.94047619047619
$ws.Range("P25").Value = 42
$ws.Range("Q25").Value = 4
$ws.Range("R25").Value = 4
$ws.Range("S25").Value = 4
$ws.Range("T25").Value = 4
$ws.Range("U25").Value = 4
$ws.Range("V25").Value = 4

$ws.Range("L26").Value = "stimuli/img_iudc4.png"
$ws.Range("M26").Value = 73.625
$ws.Range("N26").Value = 52.275
$ws.Range("O26").Value = 62.95
$ws.Range("P26").Value = 40
$ws.Range("Q26").Value = 6
$ws.Range("R26").Value = 6
$ws.Range("S26").Value = 6
$ws.Range("T26").Value = 6
$ws.Range("U26").Value = 6
$ws.Range("V26").Value = 6

$ws.Range("L27").Value = "stimuli/img_bbs77.png"
$ws.Range("M27").Value = 31.64444444444445
$ws.Range("N27").Value = 21.26666666666667
$ws.Range("O27").Value = 26.45555555555556
$ws.Range("P27").Value = 45
$ws.Range("Q27").Value = 2
$ws.Range("R27").Value = 2
$ws.Range("S27").Value = 2
$ws.Range("T27").Value = 2
$ws.Range("U27").Value = 2
$ws.Range("V27").Value = 2

$ws.Range("L28").Value = "stimuli/img_lzz3x.png"
$ws.Range("M28").Value = 18.46341463414634
$ws.Range("N28").Value = 11.92682926829268
$ws.Range("O28").Value = 15.19512195121951
$ws.Range("P28").Value = 41
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = 1
$ws.Range("S28").Value = 1
$ws.Range("T28").Value = 1
$ws.Range("U28").Value = 1
$ws.Range("V28").Value = 1

$ws.Range("L29").Value = "stimuli/img_196rk.png"
$ws.Range("M29").Value = 86.53488372093024
$ws.Range("N29").Value = 69.46511627906976
$ws.Range("O29").Value = 78
$ws.Range("P29").Value = 43
$ws.Range("Q29").Value = 9
$ws.Range("R29").Value = 9
$ws.Range("S29").Value = 9
$ws.Range("T29").Value = 9
$ws.Range("U29").Value = 9
$ws.Range("V29").Value = 9

$ws.Range("L31").Value = "stimuli/img_37hgm.png"
$ws.Range("M31").Value = 70.95454545454545
$ws.Range("N31").Value = 54.77272727272727
$ws.Range("O31").Value = 62.86363636363636
$ws.Range("P31").Value = 44
$ws.Range("Q31").Value = 6
$ws.Range("R31").Value = 6
$ws.Range("S31").Value = 6
$ws.Range("T31").Value = 6
$ws.Range("U31").Value = 6
$ws.Range("V31").Value = 6

$ws.Range("L32").Value = "stimuli/img_tbs4n.png"
$ws.Range("M32").Value = 78.95744680851064
$ws.Range("N32").Value = 58.97872340425532
$ws.Range("O32").Value = 68.96808510638297
$ws.Range("P32").Value = 47
$ws.Range("Q32").Value = 7
$ws.Range("R32").Value = 7
$ws.Range("S32").Value = 7
$ws.Range("T32").Value = 7
$ws.Range("U32").Value = 7
$ws.Range("V32").Value = 7

$ws.Range("L33").Value = "stimuli/img_rg4in.png"
$ws.Range("M33").Value = 49.3695652173913
$ws.Range("N33").Value = 30.21739130434782
$ws.Range("O33").Value = 39.79347826086956
$ws.Range("P33").Value = 46
$ws.Range("Q33").Value = 3
$ws.Range("R33").Value = 3
$ws.Range("S33").Value = 3
$ws.Range("T33").Value = 3
$ws.Range("U33").Value = 3
$ws.Range("V33").Value = 3

$ws.Range("L36").Value = "stimuli/img_eiu3c.png"
$ws.Range("M36").Value = 65.1590909090909
$ws.Range("N36").Value = 46.22727272727273
$ws.Range("O36").Value = 55.69318181818181
$ws.Range("P36").Value = 44
$ws.Range("Q36").Value = 5
$ws.Range("R36").Value = 5
$ws.Range("S36").Value = 5
$ws.Range("T36").Value = 5
$ws.Range("U36").Value = 5
$ws.Range("V36").Value = 5

$ws.Range("L41").Value = "stimuli/img_bj99b.png"
$ws.Range("M41").Value = 82.79069767441861
$ws.Range("N41").Value = 65.46511627906976
$ws.Range("O41").Value = 74.12790697674419
$ws.Range("P41").Value = 43
$ws.Range("Q41").Value = 8
$ws.Range("R41").Value = 8
$ws.Range("S41").Value = 8
$ws.Range("T41").Value = 8
$ws.Range("U41").Value = 8
$ws.Range("V41").Value = 8

